# Update the guard schedule worksheet with the new shift assignments.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Kate North)
$ws.Range("B6").Value = "OFF"
$ws.Range("F6").Value = "1:00-6:00"
$ws.Range("G6").Value = "10:30-3:30"

# Row 8 (Avery Larsen)
$ws.Range("B8").Value = "10:15-3:30"
$ws.Range("F8").Value = "10:00-3:30"
$ws.Range("G8").Value = "1:00-6:00"

# Row 9 (Austin Page)
$ws.Range("F9").Value = "3:30-8"
$ws.Range("G9").Value = "OFF"

# Row 11 (Robert Wade)
$ws.Range("B11").Value = "OFF"
$ws.Range("C11").Value = "3:30-8"

# Row 12 (Tatum Plunk)
$ws.Range("B12").Value = "OFF"

# Row 13 (Michael Vangruber)
$ws.Range("B13").Value = "OFF"
$ws.Range("G13").Value = "OFF"

# Row 14 (Jackson Blakely)
$ws.Range("B14").Value = "OFF"
$ws.Range("F14").Value = "3:30-8"

# Row 15 (Addison Clark)
$ws.Range("B15").Value = "10:15-3:30"
$ws.Range("F15").Value = "10:00-3:30"
$ws.Range("G15").Value = "OFF"

# Row 16 (Madison Johnson)
$ws.Range("B16").Value = "3:30-8"
$ws.Range("C16").Value = "10:30-3:30"
$ws.Range("E16").Value = "4:00-9"
$ws.Range("F16").Value = "10:15-3:30"
$ws.Range("G16").Value = "3:30-8"

# Row 20 (Ethan Van Horn)
$ws.Range("B20").Value = "10:30-3:30"

# Row 21 (Kai King)
$ws.Range("B21").Value = "3:30-8"
$ws.Range("E21").Value = "1:00-6:00"
$ws.Range("G21").Value = "3:30-8"

# Row 22 (Madeline Ellison)
$ws.Range("B22").Value = "10:30-3:30"
$ws.Range("F22").Value = "3:30-8"

# Row 23 (Tyler Carpenter)
$ws.Range("B23").Value = "3:30-8"
$ws.Range("C23").Value = "3:30-8"
$ws.Range("F23").Value = "OFF"

# Row 25 (Jayden Garcia)
$ws.Range("C25").Value = "OFF"
$ws.Range("E25").Value = "4:00-9"
$ws.Range("F25").Value = "3:30-8"

# Row 26 (Naya Okonkwo)
$ws.Range("B26").Value = "10:30-3:30"
$ws.Range("E26").Value = "10:30-4"

# Row 27 (Bella Hamilton)
$ws.Range("G27").Value = "10:45-3:30"

# Row 28 (Phillip Thompson)
$ws.Range("E28").Value = "10:30-4"
$ws.Range("G28").Value = "3:30-8"
